$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 <-> Row 5 got swapped (Company/Price/Day High in columns B:D).
# Use copy/paste through a scratch range so the original text formatting
# (values stored as text, e.g. "1066.00") is preserved instead of being
# reinterpreted as numbers.
$ws.Range("B4:D4").Copy()
$ws.Range("Z1:AB1").PasteSpecial()
$ws.Range("B5:D5").Copy()
$ws.Range("B4:D4").PasteSpecial()
$ws.Range("Z1:AB1").Copy()
$ws.Range("B5:D5").PasteSpecial()
$ws.Range("Z1:AB1").Clear()

# Row 10 <-> Row 11 got swapped (Company/Price/Day High in columns B:D).
$ws.Range("B10:D10").Copy()
$ws.Range("Z1:AB1").PasteSpecial()
$ws.Range("B11:D11").Copy()
$ws.Range("B10:D10").PasteSpecial()
$ws.Range("Z1:AB1").Copy()
$ws.Range("B11:D11").PasteSpecial()
$ws.Range("Z1:AB1").Clear()

$excel.CutCopyMode = $false
